$d = $word.ActiveDocument

# Locate the run of text that needs to be split into three runs:
#   "do Cục CSQLHC về TTXH cấp"  ->  "do " + "Cục CSQLHC về TTXH" + " cấp"
$searchText = "do Cục CSQLHC về TTXH cấp"
$r = $d.Content
$r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($r.Find.Found) {
    $start = $r.Start
    $end = $r.End

    $firstPart = "do "
    $middlePart = "Cục CSQLHC về TTXH"

    $len1 = $firstPart.Length
    $len2 = $middlePart.Length

    # The middle segment is re-stamped with its own (unchanged) font color.
    # Toggling the value forces Word to materialize a run boundary around
    # this sub-range, splitting the original single run into three runs
    # that share the same formatting.
    $middleRange = $d.Range($start + $len1, $start + $len1 + $len2)
    $originalColor = $middleRange.Font.Color
    $middleRange.Font.Color = 1
    $middleRange.Font.Color = $originalColor
}
